$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("DLC_List")
$ws.Activate()
$ws.Range("H2:H12").Value = "600"
$ws.Range("H2:H12").Select()
